# Locate the target paragraph: "Ajouter bloquage après clic sur inscription fosse."
$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*bloquage*" -and $p.Range.Text -like "*inscription fosse*") {
        $target = $p
        break
    }
}

$targetRange = $target.Range

# Insert a brand-new list paragraph right after the target one, before applying
# the strike-through formatting, so the new paragraph does not inherit it.
$targetRange.InsertParagraphAfter()

# The newly created paragraph is the one immediately following the target.
$newPara = $target.Next()
$newPara.Range.Text = "Dupliquer les exercices d’une palanquée d’une autre séance"

# Now strike through the whole original paragraph (its paragraph mark and all runs).
$targetRange.Font.StrikeThrough = 1
